$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the huge error text in C3 with the short placeholder text
# (matches the value already present in C4). A leading apostrophe is
# Excel's "treat as text" prefix character, so it must be doubled when
# assigned through the object model to make the apostrophe itself part
# of the stored string. Setting the value also flips the cell's style to
# a quote-prefix style, so restore the default "Normal" style afterward
# to keep the cell unstyled like the original.
$ws.Range("C3").Value2 = "''product_description'"
$ws.Range("C3").Style = "Normal"

# Append a new row 5 with the scraped data for mkdistro.com
$ws.Range("A5").Value2 = "https://www.mkdistro.com/"
$ws.Range("B5").Value2 = "ERROR"
$ws.Range("C5").Value2 = "''product_description'"
$ws.Range("C5").Style = "Normal"
